{"js": "// Insert a new \"List Bullet\" paragraph with the docente's name right\n// after the \"Docente(s) Respons\u00e1vel(eis)\" heading paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst heading = paragraphs.items.find(\n  (p) => p.text.trim() === \"Docente(s) Respons\u00e1vel(eis)\"\n);\n\nif (!heading) {\n  throw new Error('Paragraph \"Docente(s) Respons\u00e1vel(eis)\" not found.');\n}\n\nconst newParagraph = heading.insertParagraph(\n  \"5817650 - \u00c9rica Leonor Rom\u00e3o\",\n  \"After\"\n);\nnewParagraph.style = \"List Bullet\";\n\nawait context.sync();\n", "ps1": "# Insert a new \"List Bullet\" paragraph with the docente's name right\n# after the \"Docente(s) Respons\u00e1vel(eis)\" heading paragraph.\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"Docente(s) Respons\u00e1vel(eis)\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Paragraph 'Docente(s) Respons\u00e1vel(eis)' not found.\"\n}\n\n$target.Range.InsertParagraphAfter()\n$newPara = $target.Next()\n$newPara.Range.Text = \"5817650 - \u00c9rica Leonor Rom\u00e3o\"\n$newPara.Style = \"List Bullet\"\n"}
